$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.084.61"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.877.32"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5039"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3828"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08481"
$ws.Range("E9").Value = "  -8.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.245"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "1.881.65"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.56"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.217"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001099"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.15"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06668"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.03"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.084"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "28.110.02"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.583"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "2.100.44"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.65"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.60"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.51"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1047"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.630"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.620"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.704"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06548"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2174"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.223"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6531"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.247"
$ws.Range("E41").Value = "  -8.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.34"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.909"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6224"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.299"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.685"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.020"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.222"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.11"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.32"
$ws.Range("E51").Value = "  +1.64%  "
